$data = @{
    35 = @(44413, 500, 8000, 8000, 8000, 1333)
    36 = @(44377, 500, 7500, 7500, 7500, 1250)
    37 = @(44400, 500, 8000, 8000, 8000, 1333)
    38 = @(44295, 500, 8000, 8000, 8000, 1333)
    39 = @(44390, 800, 8000, 8000, 8000, 1333)
    40 = @(44298, 300, 8000, 8000, 8000, 1333)
    41 = @(44383, 600, 8000, 8000, 8000, 1333)
    42 = @(44354, 500, 7500, 7500, 7500, 1250)
    43 = @(44410, 700, 8000, 8000, 8000, 1333)
    44 = @(44412, 500, 8000, 8000, 8000, 1333)
    45 = @(44336, 600, 7000, 7000, 7000, 1167)
    46 = @(44300, 500, 8000, 8000, 8000, 1333)
    47 = @(44371, 500, 7500, 7500, 7500, 1250)
    48 = @(44320, 500, 7500, 7500, 7500, 1250)
    49 = @(44314, 300, 8000, 8000, 8000, 1333)
    50 = @(44364, 500, 7000, 7000, 7000, 1167)
    51 = @(44405, 500, 8000, 8000, 8000, 1333)
    52 = @(44327, 600, 7000, 7000, 7000, 1167)
    53 = @(44313, 500, 8000, 8000, 8000, 1333)
    54 = @(44330, 500, 7000, 7000, 7000, 1167)
    55 = @(44391, 500, 8000, 8000, 8000, 1333)
    56 = @(44351, 500, 7000, 7000, 7000, 1167)
    57 = @(44350, 500, 7000, 7000, 7000, 1167)
    58 = @(44315, 400, 8000, 8000, 8000, 1333)
    59 = @(44358, 500, 7500, 7500, 7500, 1250)
    60 = @(44389, 500, 8000, 8000, 8000, 1333)
    61 = @(44399, 600, 8000, 8000, 8000, 1333)
    62 = @(44292, 500, 8000, 8000, 8000, 1333)
    63 = @(44305, 500, 8000, 8000, 8000, 1333)
    64 = @(44294, 500, 8000, 8000, 8000, 1333)
    65 = @(44369, 500, 7000, 7000, 7000, 1167)
    66 = @(44307, 400, 8000, 8000, 8000, 1333)
    67 = @(44333, 500, 7000, 7000, 7000, 1167)
    68 = @(44309, 500, 8000, 8000, 8000, 1333)
    69 = @(44319, 500, 8000, 8000, 8000, 1333)
    70 = @(44316, 500, 8000, 8000, 8000, 1333)
    71 = @(44301, 500, 8000, 8000, 8000, 1333)
    72 = @(44370, 500, 7500, 7500, 7500, 1250)
    73 = @(44326, 500, 7000, 7000, 7000, 1167)
    74 = @(44382, 500, 8000, 8000, 8000, 1333)
    75 = @(44398, 500, 8000, 8000, 8000, 1333)
    76 = @(44355, 500, 7500, 7500, 7500, 1250)
    77 = @(44343, 500, 7000, 7000, 7000, 1167)
    78 = @(44376, 500, 7500, 7500, 7500, 1250)
    79 = @(44334, 500, 7000, 7000, 7000, 1167)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("J$r").Value = $vals[1]
    $ws.Range("K$r").Value = $vals[2]
    $ws.Range("L$r").Value = $vals[3]
    $ws.Range("M$r").Value = $vals[4]
    $ws.Range("P$r").Value = $vals[5]
}

# Row 79 is brand new: fill the static columns by copying row 78's values
# (use .Value2 for reads - .Value getter is unreliable in this host)
$ws.Range("A79").Value = $ws.Range("A78").Value2
$ws.Range("B79").Value = $ws.Range("B78").Value2
$ws.Range("C79").Value = $ws.Range("C78").Value2
$ws.Range("E79").Value = $ws.Range("E78").Value2
$ws.Range("F79").Value = $ws.Range("F78").Value2
$ws.Range("G79").Value = $ws.Range("G78").Value2
$ws.Range("H79").Value = $ws.Range("H78").Value2
$ws.Range("I79").Value = $ws.Range("I78").Value2
$ws.Range("N79").Value = $ws.Range("N78").Value2
$ws.Range("O79").Value = $ws.Range("O78").Value2
$ws.Range("Q79").Value = $ws.Range("Q78").Value2
$ws.Range("R79").Value = $ws.Range("R78").Value2

# Match the date-style formatting used by the rest of column D
$ws.Range("D79").NumberFormat = "YYYY-MM-DD HH:MM:SS"
